$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line entries (line7, line8) were inserted into the shared-string
# table right after line6 / before extr1. This shifts the meaning of the
# existing "name" cells in rows 8-15 (their stored shared-string index did
# not change, but since two new strings were spliced in earlier in the
# table, those rows now display the next names in sequence). Re-assert the
# display text for every name cell from row 8 onward so the workbook ends
# up with the correct strings in the correct cells.
$ws.Range("B8").Value  = "line7"
$ws.Range("B9").Value  = "line8"
$ws.Range("B10").Value = "extr1"
$ws.Range("B11").Value = "extr2"
$ws.Range("B12").Value = "extr3"
$ws.Range("B13").Value = "extr4"
$ws.Range("B14").Value = "extr5"
$ws.Range("B15").Value = "extr6"

# Updated data values for existing rows 8-14
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("C12").Value = 10

$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# Row 15 now holds what used to be the extr8 data slot, with new values
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# Grow the table by two new rows (16 and 17), matching the formatting of
# the last existing data row (15).
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A17:E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
